$d = $word.ActiveDocument

$newText = "Du deltar i en världsomspännande kampanj för att observera och rapportera de svagaste synliga stjärnorna, som ett mått på ljusföroreningarna på orten. Genom att hitta och observera Lejonets konstellation på natthimlen kan folk i hela världen lära sig hur belysningen i våra samhällen och omgivningar bidrar till ljusföroreningar. Era bidrag till online-databasen hjälper till att dokumentera den synliga natthimlens över hela världen."

# Collect the paragraphs whose whole text is one of the two legacy blurbs
# ("Kampanjdatum för Perseus 2018: ..." or the broken-up "Du deltar ...")
# so we can replace each of them, in full, with a single plain run
# containing the new consolidated sentence.
$targets = New-Object System.Collections.ArrayList
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Kampanjdatum för*") {
        [void]$targets.Add($p.Range)
    }
    elseif ($t -like "Du delta*") {
        [void]$targets.Add($p.Range)
    }
}

foreach ($rng in $targets) {
    $rng.End = $rng.End - 1
    $rng.Delete()
    $rng.InsertAfter($newText)
}
